$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OVERVIEW")

# Row 46 (Prague 9): status "Finishing" -> "Done"
$ws.Range("E46").Value = "Done"

# Row 4 (Las Huelgas Codex): "beginning with pieces ascribed to Vitry" -> "finished pieces ascribe to Vitry";
# new "partly" note in C4
$ws.Range("F4").Value = "finished pieces ascribe to Vitry"
$ws.Range("C4").Value = "partly"

# Row 9 (1 Fauvel): new "partly" note in C9 and new priority note in E9
$ws.Range("C9").Value = "partly"
$ws.Range("E9").Value = "process high priority, transcribe low priority"

# Row 39 (Du Fay works entry above "Jeremy Jennings" row 40): add "Corwyn Wyatt" assignment
$ws.Range("B39").Value = "Corwyn Wyatt"

# Update the active cell/selection on the sheet to E10
$ws.Activate() | Out-Null
$ws.Range("E10").Select() | Out-Null
